# Fruta / hortaliza, semanal
# The data rows (2-6) get re-sorted/rotated: for columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg), the content of each row is replaced by the content that used
# to live in a different row, per the mapping below (row -> source row):
#   2 <- 3, 3 <- 5, 4 <- 6, 5 <- 2, 6 <- 4
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are identical across all data rows
# in this sheet already, so they do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot current values of every affected cell (rows 2-6) before writing
# anything, so the row-to-row copy below never reads an already-overwritten
# cell.
$snapshot = @{}
foreach ($r in 2..6) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# row -> source row it should take its values from
$mapping = @{ 2 = 3; 3 = 5; 4 = 6; 5 = 2; 6 = 4 }

foreach ($destRow in 2..6) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot["$srcRow-$c"]
    }
}
